# Add a new slide ("Static method") at the end of the deck, using the
# "Title and Content" layout (CustomLayout #2 / slideLayout2.xml), matching
# the pattern already used by other "code sample" slides in this deck
# (e.g. slide2 which also uses Title + idx=1 content placeholder).

$p = $ppt.ActivePresentation

$layoutIndex = 2   # "Title and Content" custom layout
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, $layoutIndex)

# ---- Title -----------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Static method"

# ---- Body content ------------------------------------------------------
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# Paragraph 1
$tr.Text = "Để định nghĩa method dạng static thì thêm từ khóa"

# Paragraph 2: "@staticmethod" (bold, italic, red), no bullet
$null = $tr.InsertAfter("`r@staticmethod")

# Paragraph 3
$null = $tr.InsertAfter("`rVí dụ:")

# Paragraph 4: class class_name:
$null = $tr.InsertAfter("`rclass class_name:")

# Paragraph 5: TAB + @staticmethod (italic, red), no bullet
$null = $tr.InsertAfter("`r`t@staticmethod")

# Paragraph 6: TAB def static_method_name(param_list):
$null = $tr.InsertAfter("`r`tdef static_method_name(param_list):")

# Paragraph 7: TAB TAB pass
$null = $tr.InsertAfter("`r`t`tpass")

# Paragraph 8
$null = $tr.InsertAfter("`rCách gọi:")

# Paragraph 9: class_name.static_method_name()
$null = $tr.InsertAfter("`rclass_name.static_method_name()")

# Uniform base size for the whole placeholder
$tr.Font.Size = 28

$paras = $tr.Paragraphs()

# --- Paragraph 1: "Để định nghĩa method dạng static thì thêm từ khóa"
# (plain run, size 28 only - already set above)

# --- Paragraph 2: "@staticmethod" -> bold, italic, red; no bullet
$para2 = $tr.Paragraphs(2,1)
$para2.ParagraphFormat.Bullet.Visible = $false
$para2.Font.Bold = $true
$para2.Font.Italic = $true
$para2.Font.Color.RGB = 255

# --- Paragraph 3: "Ví dụ:" (plain)

# --- Paragraph 4: "class class_name:" -> "class_name" bold + green
$para4 = $tr.Paragraphs(4,1)
$cname4 = $para4.Characters(7, 10)
$cname4.Font.Bold = $true
$cname4.Font.Color.RGB = 32768

# --- Paragraph 5: "\t@staticmethod" -> italic, red; no bullet
$para5 = $tr.Paragraphs(5,1)
$para5.ParagraphFormat.Bullet.Visible = $false
$atmethod5 = $para5.Characters(2, 13)
$atmethod5.Font.Italic = $true
$atmethod5.Font.Color.RGB = 255

# --- Paragraph 6: "\tdef static_method_name(param_list):" -> method name bold + maroon; no bullet
$para6 = $tr.Paragraphs(6,1)
$para6.ParagraphFormat.Bullet.Visible = $false
$mname6 = $para6.Characters(6, 18)
$mname6.Font.Bold = $true
$mname6.Font.Color.RGB = 128

# --- Paragraph 7: "\t\tpass" -> no bullet
$para7 = $tr.Paragraphs(7,1)
$para7.ParagraphFormat.Bullet.Visible = $false

# --- Paragraph 8: "Cách gọi:" (plain)

# --- Paragraph 9: "class_name.static_method_name()" -> level 2, no bullet,
#                  "class_name" bold+green, "static_method_name" bold+maroon
$para9 = $tr.Paragraphs(9,1)
$para9.IndentLevel = 2
$para9.ParagraphFormat.Bullet.Visible = $false
$cname9 = $para9.Characters(1, 10)
$cname9.Font.Bold = $true
$cname9.Font.Color.RGB = 32768
$mname9 = $para9.Characters(12, 18)
$mname9.Font.Bold = $true
$mname9.Font.Color.RGB = 128

Write-Host ("Slides=" + $p.Slides.Count)
